$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (class 0.0)
$ws.Range("B2").Value = 0.8976377952755905
$ws.Range("C2").Value = 0.8837209302325582
$ws.Range("D2").Value = 0.8906249999999999

# Row 3 (class 1.0)
$ws.Range("B3").Value = 0.8728813559322034
$ws.Range("D3").Value = 0.8803418803418803

# Row 4 (accuracy)
$ws.Range("B4").Value = 0.8857142857142857
$ws.Range("C4").Value = 0.8857142857142857
$ws.Range("D4").Value = 0.8857142857142857
$ws.Range("E4").Value = 0.8857142857142857

# Row 5 (macro avg)
$ws.Range("B5").Value = 0.885259575603897
$ws.Range("C5").Value = 0.8858259823576584
$ws.Range("D5").Value = 0.8854834401709402

# Row 6 (weighted avg)
$ws.Range("B6").Value = 0.8859163790966806
$ws.Range("C6").Value = 0.8857142857142857
$ws.Range("D6").Value = 0.8857562576312575
